$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddCustomerTest")

$ws.Range("A2").Value = "Dave"
$ws.Range("B2").Value = "Smith"

$ws.Range("A3").Value = "Martin"
$ws.Range("B3").Value = "Byrne"

$ws.Range("A4").Value = "Alana"
$ws.Range("B4").Value = "Curran"

$ws.Range("A5").Value = "Paul"
$ws.Range("B5").Value = "Jones"

$ws.Range("B5").Select()
